# Add IPO dates for airlines
# Row 26 was a blank spacer row under the "Profile" block (B24=Engine Type,
# B25=Founded). Fill it in with a new "IPO" label/value pair, matching the
# formatting already used by the rows above it (B24/B25, C24/C25).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B26").Value = "IPO"
$ws.Range("C26").Value = 2015

# Match the center-aligned style already used by C24:D24 / C25:D25.
$ws.Range("C26:D26").HorizontalAlignment = -4108
